# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45175 to serial date 45177 (2023-09-06 -> 2023-09-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row on the sheet.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)  # Column C
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
